$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updating team specific time data (transition matrix probabilities)
# values recomputed with new underlying data
$ws.Range("B2").Value = 0.1161971830985915
$ws.Range("C2").Value = 0.6338028169014085
$ws.Range("J2").Value = 0.0176056338028169
$ws.Range("P2").Value = 0.1443661971830986
$ws.Range("S2").Value = 0.0880281690140845
$ws.Range("B3").Value = 0.0101010101010101
$ws.Range("C3").Value = 0.0202020202020202
$ws.Range("J3").Value = 0.04545454545454546
$ws.Range("P3").Value = 0.7525252525252525
$ws.Range("S3").Value = 0.1717171717171717
$ws.Range("J4").Value = 0.1428571428571428
$ws.Range("P4").Value = 0.6904761904761905
$ws.Range("S4").Value = 0.1666666666666667
$ws.Range("B6").Value = 0.06756756756756757
$ws.Range("D6").Value = 0.01801801801801802
$ws.Range("F6").Value = 0.05405405405405406
$ws.Range("J6").Value = 0.2702702702702703
$ws.Range("O6").Value = 0.01801801801801802
$ws.Range("Q6").Value = 0.1576576576576577
$ws.Range("R6").Value = 0.04054054054054054
$ws.Range("S6").Value = 0.3738738738738739
$ws.Range("B7").Value = 0.1055276381909548
$ws.Range("D7").Value = 0.005025125628140704
$ws.Range("F7").Value = 0.04522613065326633
$ws.Range("J7").Value = 0.07035175879396985
$ws.Range("O7").Value = 0.1055276381909548
$ws.Range("Q7").Value = 0.1608040201005025
$ws.Range("R7").Value = 0.1055276381909548
$ws.Range("S7").Value = 0.4020100502512563
$ws.Range("B8").Value = 0.1337448559670782
$ws.Range("D8").Value = 0.0205761316872428
$ws.Range("F8").Value = 0.06790123456790123
$ws.Range("J8").Value = 0.1172839506172839
$ws.Range("O8").Value = 0.01646090534979424
$ws.Range("Q8").Value = 0.1831275720164609
$ws.Range("R8").Value = 0.0720164609053498
$ws.Range("S8").Value = 0.3888888888888889
$ws.Range("B9").Value = 0.09625668449197861
$ws.Range("D9").Value = 0.0213903743315508
$ws.Range("F9").Value = 0.0748663101604278
$ws.Range("J9").Value = 0.08021390374331551
$ws.Range("O9").Value = 0.0213903743315508
$ws.Range("Q9").Value = 0.1925133689839572
$ws.Range("R9").Value = 0.05882352941176471
$ws.Range("S9").Value = 0.4545454545454545
$ws.Range("B10").Value = 0.08972392638036809
$ws.Range("D10").Value = 0.01917177914110429
$ws.Range("E10").Value = 0.001533742331288344
$ws.Range("F10").Value = 0.06441717791411043
$ws.Range("J10").Value = 0.1748466257668712
$ws.Range("O10").Value = 0.02530674846625767
$ws.Range("Q10").Value = 0.196319018404908
$ws.Range("R10").Value = 0.08282208588957055
$ws.Range("S10").Value = 0.3458588957055215
$ws.Range("G11").Value = 0.1538461538461539
$ws.Range("J11").Value = 0.09523809523809523
$ws.Range("K11").Value = 0.1428571428571428
$ws.Range("L11").Value = 0.608058608058608
$ws.Range("G12").Value = 0.7261904761904762
$ws.Range("J12").Value = 0.2142857142857143
$ws.Range("K12").Value = 0.005952380952380952
$ws.Range("L12").Value = 0.02976190476190476
$ws.Range("S12").Value = 0.02380952380952381
$ws.Range("F13").Value = 0.0196078431372549
$ws.Range("G13").Value = 0.7058823529411765
$ws.Range("J13").Value = 0.2352941176470588
$ws.Range("S13").Value = 0.0392156862745098
$ws.Range("F15").Value = 0.01904761904761905
$ws.Range("H15").Value = 0.1523809523809524
$ws.Range("I15").Value = 0.04761904761904762
$ws.Range("J15").Value = 0.3523809523809524
$ws.Range("K15").Value = 0.07142857142857142
$ws.Range("M15").Value = 0.01428571428571429
$ws.Range("O15").Value = 0.04761904761904762
$ws.Range("S15").Value = 0.2952380952380952
$ws.Range("F16").Value = 0.0892018779342723
$ws.Range("H16").Value = 0.2065727699530517
$ws.Range("I16").Value = 0.07981220657276995
$ws.Range("J16").Value = 0.3661971830985916
$ws.Range("K16").Value = 0.08450704225352113
$ws.Range("M16").Value = 0.02347417840375587
$ws.Range("O16").Value = 0.02816901408450704
$ws.Range("S16").Value = 0.1220657276995305
$ws.Range("F17").Value = 0.01111111111111111
$ws.Range("H17").Value = 0.1755555555555555
$ws.Range("I17").Value = 0.1088888888888889
$ws.Range("J17").Value = 0.3977777777777778
$ws.Range("K17").Value = 0.08888888888888889
$ws.Range("M17").Value = 0.02666666666666667
$ws.Range("O17").Value = 0.06
$ws.Range("S17").Value = 0.1311111111111111
$ws.Range("F18").Value = 0.02162162162162162
$ws.Range("H18").Value = 0.1567567567567568
$ws.Range("I18").Value = 0.07567567567567568
$ws.Range("J18").Value = 0.4216216216216216
$ws.Range("K18").Value = 0.1405405405405405
$ws.Range("M18").Value = 0.03243243243243243
$ws.Range("O18").Value = 0.04864864864864865
$ws.Range("S18").Value = 0.1027027027027027
$ws.Range("F19").Value = 0.01193317422434368
$ws.Range("H19").Value = 0.233890214797136
$ws.Range("I19").Value = 0.07875894988066826
$ws.Range("J19").Value = 0.3707239459029435
$ws.Range("K19").Value = 0.1081941129673827
$ws.Range("M19").Value = 0.02068416865552904
$ws.Range("N19").Value = 0.0007955449482895784
$ws.Range("O19").Value = 0.05727923627684964
$ws.Range("S19").Value = 0.1177406523468576
